# Fruta / hortaliza, semanal
# Updates rows 23-28 with new weekly data and appends three new rows (29-31)
# that carry forward the prior week's values that were previously in rows 26-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: Dina/Especial -> Castle Brite/Especial, new week, new region ---
$ws.Range("D23").Value = 44543
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("M23").Value = 300
$ws.Range("R23").Value = "Región de O'Higgins"

# --- Row 24: Castle Brite/Primera, new week, new region, new prices ---
$ws.Range("D24").Value = 44543
$ws.Range("M24").Value = 400
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 20500
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 1139

# --- Row 25: Castle Brite/Segunda, new week, new region, new prices ---
$ws.Range("D25").Value = 44543
$ws.Range("M25").Value = 400
$ws.Range("N25").Value = 15000
$ws.Range("O25").Value = 16000
$ws.Range("P25").Value = 15500
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 861

# --- Row 26: now Dina/Especial (same data that used to sit in row 23) ---
$ws.Range("D26").Value = 44187
$ws.Range("K26").Value = "Dina"
$ws.Range("M26").Value = 240
$ws.Range("N26").Value = 22000
$ws.Range("O26").Value = 23000
$ws.Range("P26").Value = 22500
$ws.Range("Q26").Value = "$/caja 18 kilos"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 1250
$ws.Range("T26").Value = 18

# --- Row 27: now Castle Brite/Primera (same data that used to sit in row 24) ---
$ws.Range("D27").Value = 44536
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 22000
$ws.Range("O27").Value = 23000
$ws.Range("P27").Value = 22500
$ws.Range("Q27").Value = "$/caja 18 kilos"
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 1250
$ws.Range("T27").Value = 18

# --- Row 28: now Castle Brite/Segunda (same data that used to sit in row 25) ---
$ws.Range("D28").Value = 44536
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 19000
$ws.Range("P28").Value = 18500
$ws.Range("Q28").Value = "$/caja 18 kilos"
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 1028
$ws.Range("T28").Value = 18

# --- New row 29: Castle Brite/Especial (the data that used to sit in row 26) ---
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44540
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103003
$ws.Range("J29").Value = "Damasco"
$ws.Range("K29").Value = "Castle Brite"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 140
$ws.Range("N29").Value = 24500
$ws.Range("O29").Value = 25000
$ws.Range("P29").Value = 24750
$ws.Range("Q29").Value = "$/caja 15 kilos"
$ws.Range("R29").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S29").Value = 1650
$ws.Range("T29").Value = 15

# --- New row 30: Castle Brite/Primera (the data that used to sit in row 27) ---
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44540
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103003
$ws.Range("J30").Value = "Damasco"
$ws.Range("K30").Value = "Castle Brite"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 22500
$ws.Range("O30").Value = 23000
$ws.Range("P30").Value = 22750
$ws.Range("Q30").Value = "$/caja 15 kilos"
$ws.Range("R30").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S30").Value = 1517
$ws.Range("T30").Value = 15

# --- New row 31: Castle Brite/Segunda (the data that used to sit in row 28) ---
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44540
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103003
$ws.Range("J31").Value = "Damasco"
$ws.Range("K31").Value = "Castle Brite"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 18000
$ws.Range("O31").Value = 18500
$ws.Range("P31").Value = 18250
$ws.Range("Q31").Value = "$/caja 15 kilos"
$ws.Range("R31").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S31").Value = 1217
$ws.Range("T31").Value = 15

# Ensure date cells carry the workbook's existing date style (style index 2,
# the same custom date format already used by D2:D28).
$ws.Range("D29").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("D30").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("D31").NumberFormat = $ws.Range("D28").NumberFormat
